# Add a new "Screen Resolution" column (N) and move "Computer Company" to the
# new last column (O). Populate the new Screen Resolution column with "FHD"
# for every data row, and normalize some "Graphic Card" (M) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before N ("Computer Company") so it shifts to O,
#    then give the freed-up N column the header "Screen Resolution".
#    Inserting the column shifts formatting along with the data, so N1
#    already carries the header look (bold + border); just set its text.
$ws.Columns.Item(14).Insert()
$ws.Cells.Item(1, 14).Value = "Screen Resolution"

# 2. Fill column N (Screen Resolution) with "FHD" for every data row (2-65).
for ($r = 2; $r -le 65; $r++) {
    $ws.Cells.Item($r, 14).Value = "FHD"
}

# 3. Normalize specific "Graphic Card" values in column M.
$graphicCardUpdates = @{
    3  = "NVIDIA® GTX"
    4  = "NVIDIA® GeForce® GTX"
    5  = "NVIDIA® GeForce®"
    6  = "NVIDIA® GeForce® GTX"
    7  = "NVIDIA® GTX"
    8  = "NVIDIA® GeForce®"
    9  = "NVIDIA® Quadro®"
    10 = "NVIDIA® Quadro®"
    11 = "Intel® UHD Graphics"
    12 = "Intel® UHD Graphics"
    13 = "Intel® Iris® Xe Graphics"
    14 = "NVIDIA® GeForce® GTX"
    15 = "Intel® UHD Graphics"
    16 = "Intel® UHD Graphics"
    17 = "Intel® UHD Graphics"
    18 = "Intel® UHD Graphics"
    19 = "Intel® UHD Graphics"
    23 = "Intel® Iris® Xe Graphics"
    24 = "NVIDIA® GTX"
    25 = "NVIDIA® GeForce® RTX™"
    26 = "NVIDIA® GeForce® RTX™"
    27 = "NVIDIA® Quadro®"
    28 = "NVIDIA® Quadro®"
    29 = "NVIDIA® Quadro®"
    30 = "NVIDIA® Quadro®"
    31 = "NVIDIA® Quadro®"
    32 = "NVIDIA® Quadro®"
    41 = "Radeon Pro"
    42 = "Radeon Pro"
    43 = "AMD Radeon™ Graphics"
    44 = "Intel® UHD Graphics"
    45 = "Intel® UHD Graphics"
    46 = "Intel® UHD Graphics"
    47 = "NVIDIA® GeForce® RTX™"
    48 = "AMD Radeon™ Graphics"
    49 = "AMD Radeon™ Graphics"
    50 = "Intel® UHD Graphics"
    51 = "Intel® UHD Graphics"
    52 = "Intel® UHD Graphics"
    53 = "AMD Radeon™ Graphics"
    54 = "Intel® UHD Graphics"
    56 = "Radeon Pro"
    59 = "Radeon Pro"
}

foreach ($row in $graphicCardUpdates.Keys) {
    $ws.Cells.Item($row, 13).Value = $graphicCardUpdates[$row]
}
